# Weekly refresh of the "Albahaca" price sheet: a new weekly record is
# inserted at row 57 (pushing the existing rows 57-65 down to 58-66).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57; this shifts rows 57..65 down to 58..66
# and carries the existing row formatting (e.g. the date style on column D)
# down with them.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly record.
$ws.Cells.Item(57,1).Value  = 1
$ws.Cells.Item(57,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(57,4).Value  = 45106
$ws.Cells.Item(57,5).Value  = 15
$ws.Cells.Item(57,6).Value  = 100112052
$ws.Cells.Item(57,7).Value  = "Albahaca"
$ws.Cells.Item(57,8).Value  = "Sin especificar"
$ws.Cells.Item(57,9).Value  = "Primera"
$ws.Cells.Item(57,10).Value = 370
$ws.Cells.Item(57,11).Value = 1200
$ws.Cells.Item(57,12).Value = 1500
$ws.Cells.Item(57,13).Value = 1346
$ws.Cells.Item(57,14).Value = "`$/paquete"
$ws.Cells.Item(57,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(57,16).Value = 1346
$ws.Cells.Item(57,17).Value = 1
$ws.Cells.Item(57,18).Value = "Hortaliza"
